# Update the cryptocurrency market-data table (Sheet1) to the 2023-10-22 snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 29957
$ws.Cells.Item(2, 5).Value = 584808995221
$ws.Cells.Item(2, 6).Value = 10849693265
$ws.Cells.Item(2, 7).Value = 0.65991

$ws.Cells.Item(3, 4).Value = 1632.8
$ws.Cells.Item(3, 5).Value = 196053085743
$ws.Cells.Item(3, 6).Value = 8864498457
$ws.Cells.Item(3, 7).Value = 1.63451

$ws.Cells.Item(4, 4).Value = 1.001
$ws.Cells.Item(4, 5).Value = 84206193325
$ws.Cells.Item(4, 6).Value = 20825028473
$ws.Cells.Item(4, 7).Value = 0.08513999999999999

$ws.Cells.Item(5, 4).Value = 214.51
$ws.Cells.Item(5, 5).Value = 33004927207
$ws.Cells.Item(5, 6).Value = 315721272
$ws.Cells.Item(5, 7).Value = 0.72909

$ws.Cells.Item(6, 4).Value = 0.517364
$ws.Cells.Item(6, 5).Value = 27633664259
$ws.Cells.Item(6, 6).Value = 780442138
$ws.Cells.Item(6, 7).Value = 0.1232

$ws.Cells.Item(7, 4).Value = 1.001
$ws.Cells.Item(7, 5).Value = 25574900868
$ws.Cells.Item(7, 6).Value = 3849453521
$ws.Cells.Item(7, 7).Value = 0.15668

$ws.Cells.Item(8, 4).Value = 1630.58
$ws.Cells.Item(8, 5).Value = 14394544085
$ws.Cells.Item(8, 6).Value = 5090755
$ws.Cells.Item(8, 7).Value = 1.55775

$ws.Cells.Item(9, 4).Value = 28.53
$ws.Cells.Item(9, 5).Value = 11834165412
$ws.Cells.Item(9, 6).Value = 1100945229
$ws.Cells.Item(9, 7).Value = -2.58965

$ws.Cells.Item(10, 4).Value = 0.258067
$ws.Cells.Item(10, 5).Value = 9008125407
$ws.Cells.Item(10, 6).Value = 154301151
$ws.Cells.Item(10, 7).Value = 0.4562

$ws.Cells.Item(11, 4).Value = 0.060826
$ws.Cells.Item(11, 5).Value = 8599337024
$ws.Cells.Item(11, 6).Value = 257674055
$ws.Cells.Item(11, 7).Value = 0.69101

$ws.Cells.Item(12, 4).Value = 0.090682
$ws.Cells.Item(12, 5).Value = 8050993822
$ws.Cells.Item(12, 6).Value = 239012195
$ws.Cells.Item(12, 7).Value = 0.07381

$ws.Cells.Item(13, 4).Value = 2.19
$ws.Cells.Item(13, 5).Value = 7467835445
$ws.Cells.Item(13, 6).Value = 11540544
$ws.Cells.Item(13, 7).Value = 3.96587

$ws.Cells.Item(14, 2).Value = "MATIC"
$ws.Cells.Item(14, 3).Value = "Polygon"
$ws.Cells.Item(14, 4).Value = 0.562612
$ws.Cells.Item(14, 5).Value = 5226864715
$ws.Cells.Item(14, 6).Value = 288390118
$ws.Cells.Item(14, 7).Value = 1.07601

$ws.Cells.Item(15, 2).Value = "LINK"
$ws.Cells.Item(15, 3).Value = "Chainlink"
$ws.Cells.Item(15, 4).Value = 9.16
$ws.Cells.Item(15, 5).Value = 5103687268
$ws.Cells.Item(15, 6).Value = 1433618634
$ws.Cells.Item(15, 7).Value = 12.96827

$ws.Cells.Item(16, 2).Value = "DOT"
$ws.Cells.Item(16, 3).Value = "Polkadot"
$ws.Cells.Item(16, 4).Value = 3.85
$ws.Cells.Item(16, 5).Value = 4943845852
$ws.Cells.Item(16, 6).Value = 163364869
$ws.Cells.Item(16, 7).Value = 1.25749

$ws.Cells.Item(17, 4).Value = 29947
$ws.Cells.Item(17, 5).Value = 4878864838
$ws.Cells.Item(17, 6).Value = 166556288
$ws.Cells.Item(17, 7).Value = 0.70744

$ws.Cells.Item(18, 1).Value = 18
$ws.Cells.Item(18, 4).Value = 241.89
$ws.Cells.Item(18, 5).Value = 4726605181
$ws.Cells.Item(18, 6).Value = 259004248
$ws.Cells.Item(18, 7).Value = -0.13808

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "LTC"
$ws.Cells.Item(19, 3).Value = "Litecoin"
$ws.Cells.Item(19, 4).Value = 64.16
$ws.Cells.Item(19, 5).Value = 4725106929
$ws.Cells.Item(19, 6).Value = 367994420
$ws.Cells.Item(19, 7).Value = -0.18957

$ws.Cells.Item(20, 2).Value = "SHIB"
$ws.Cells.Item(20, 3).Value = "Shiba Inu"
$ws.Cells.Item(20, 4).Value = 0.000007
$ws.Cells.Item(20, 5).Value = 4124028390
$ws.Cells.Item(20, 6).Value = 105502611
$ws.Cells.Item(20, 7).Value = 0.04033

$ws.Cells.Item(21, 4).Value = 0.999839
$ws.Cells.Item(21, 5).Value = 3807741408
$ws.Cells.Item(21, 6).Value = 117659311
$ws.Cells.Item(21, 7).Value = 0.10244

$ws.Cells.Item(22, 4).Value = 3.83
$ws.Cells.Item(22, 5).Value = 3593981215
$ws.Cells.Item(22, 6).Value = 244056
$ws.Cells.Item(22, 7).Value = 0.01329

$ws.Cells.Item(23, 2).Value = "AVAX"
$ws.Cells.Item(23, 3).Value = "Avalanche"
$ws.Cells.Item(23, 4).Value = 9.75
$ws.Cells.Item(23, 5).Value = 3460709312
$ws.Cells.Item(23, 6).Value = 256660296
$ws.Cells.Item(23, 7).Value = 2.32995

$ws.Cells.Item(24, 2).Value = "TUSD"
$ws.Cells.Item(24, 3).Value = "TrueUSD"
$ws.Cells.Item(24, 4).Value = 0.9987740000000001
$ws.Cells.Item(24, 5).Value = 3363082850
$ws.Cells.Item(24, 6).Value = 216297492
$ws.Cells.Item(24, 7).Value = -0.05971

$ws.Cells.Item(25, 4).Value = 4.13
$ws.Cells.Item(25, 5).Value = 3107938437
$ws.Cells.Item(25, 6).Value = 101989244
$ws.Cells.Item(25, 7).Value = 2.31216

$ws.Cells.Item(26, 4).Value = 0.109592
$ws.Cells.Item(26, 5).Value = 3042204571
$ws.Cells.Item(26, 6).Value = 79748345
$ws.Cells.Item(26, 7).Value = 0.38027

$ws.Cells.Item(27, 4).Value = 157.82
$ws.Cells.Item(27, 5).Value = 2878984486
$ws.Cells.Item(27, 6).Value = 67192067
$ws.Cells.Item(27, 7).Value = 0.92359

$ws.Cells.Item(28, 4).Value = 43.66
$ws.Cells.Item(28, 5).Value = 2619542072
$ws.Cells.Item(28, 6).Value = 6072831
$ws.Cells.Item(28, 7).Value = -0.32057

$ws.Cells.Item(29, 4).Value = 15.5
$ws.Cells.Item(29, 5).Value = 2217780253
$ws.Cells.Item(29, 6).Value = 101375393
$ws.Cells.Item(29, 7).Value = -0.22062

$ws.Cells.Item(30, 5).Value = 2094650385
$ws.Cells.Item(30, 6).Value = 2515079753
$ws.Cells.Item(30, 7).Value = -0.03057

$ws.Cells.Item(31, 4).Value = 6.6
$ws.Cells.Item(31, 5).Value = 1931582900
$ws.Cells.Item(31, 6).Value = 124031551
$ws.Cells.Item(31, 7).Value = 1.95326

$ws.Cells.Item(32, 4).Value = 0.04868511
$ws.Cells.Item(32, 5).Value = 1628148619
$ws.Cells.Item(32, 6).Value = 40929832
$ws.Cells.Item(32, 7).Value = 1.25622

$ws.Cells.Item(33, 4).Value = 3.37
$ws.Cells.Item(33, 5).Value = 1545506907
$ws.Cells.Item(33, 6).Value = 135260877
$ws.Cells.Item(33, 7).Value = 3.91071

$ws.Cells.Item(34, 2).Value = "APT"
$ws.Cells.Item(34, 3).Value = "Aptos"
$ws.Cells.Item(34, 4).Value = 6.1
$ws.Cells.Item(34, 5).Value = 1501761414
$ws.Cells.Item(34, 6).Value = 451877110
$ws.Cells.Item(34, 7).Value = 17.29901

$ws.Cells.Item(35, 2).Value = "LDO"
$ws.Cells.Item(35, 3).Value = "Lido DAO"
$ws.Cells.Item(35, 4).Value = 1.65
$ws.Cells.Item(35, 5).Value = 1464195804
$ws.Cells.Item(35, 6).Value = 42586699
$ws.Cells.Item(35, 7).Value = 4.74214

$ws.Cells.Item(36, 2).Value = "ICP"
$ws.Cells.Item(36, 3).Value = "Internet Computer"
$ws.Cells.Item(36, 4).Value = 3.17
$ws.Cells.Item(36, 5).Value = 1416024233
$ws.Cells.Item(36, 6).Value = 17251454
$ws.Cells.Item(36, 7).Value = -0.09295

$ws.Cells.Item(37, 2).Value = "CRO"
$ws.Cells.Item(37, 3).Value = "Cronos"
$ws.Cells.Item(37, 4).Value = 0.052775
$ws.Cells.Item(37, 5).Value = 1389772735
$ws.Cells.Item(37, 6).Value = 5888898
$ws.Cells.Item(37, 7).Value = 0.31882

$ws.Cells.Item(38, 4).Value = 90.17
$ws.Cells.Item(38, 5).Value = 1311303058
$ws.Cells.Item(38, 6).Value = 20762406
$ws.Cells.Item(38, 7).Value = 3.75501

$ws.Cells.Item(39, 2).Value = "MKR"
$ws.Cells.Item(39, 3).Value = "Maker"
$ws.Cells.Item(39, 4).Value = 1426.64
$ws.Cells.Item(39, 5).Value = 1283850151
$ws.Cells.Item(39, 6).Value = 46788514
$ws.Cells.Item(39, 7).Value = -0.02245

$ws.Cells.Item(40, 2).Value = "VET"
$ws.Cells.Item(40, 3).Value = "VeChain"
$ws.Cells.Item(40, 4).Value = 0.0168811
$ws.Cells.Item(40, 5).Value = 1226940730
$ws.Cells.Item(40, 6).Value = 34526179
$ws.Cells.Item(40, 7).Value = -0.28323

$ws.Cells.Item(41, 2).Value = "MNT"
$ws.Cells.Item(41, 3).Value = "Mantle"
$ws.Cells.Item(41, 4).Value = 0.382322
$ws.Cells.Item(41, 5).Value = 1185564412
$ws.Cells.Item(41, 6).Value = 39905135
$ws.Cells.Item(41, 7).Value = 1.28823

$ws.Cells.Item(42, 2).Value = "OP"
$ws.Cells.Item(42, 3).Value = "Optimism"
$ws.Cells.Item(42, 4).Value = 1.26
$ws.Cells.Item(42, 5).Value = 1105751472
$ws.Cells.Item(42, 6).Value = 80085346
$ws.Cells.Item(42, 7).Value = 0.86395

$ws.Cells.Item(43, 2).Value = "AAVE"
$ws.Cells.Item(43, 3).Value = "Aave"
$ws.Cells.Item(43, 4).Value = 75.64
$ws.Cells.Item(43, 5).Value = 1101568659
$ws.Cells.Item(43, 6).Value = 268628968
$ws.Cells.Item(43, 7).Value = 11.96005

$ws.Cells.Item(44, 2).Value = "ARB"
$ws.Cells.Item(44, 3).Value = "Arbitrum"
$ws.Cells.Item(44, 4).Value = 0.826438
$ws.Cells.Item(44, 5).Value = 1053330801
$ws.Cells.Item(44, 6).Value = 154914980
$ws.Cells.Item(44, 7).Value = 0.54126

$ws.Cells.Item(45, 2).Value = "KAS"
$ws.Cells.Item(45, 3).Value = "Kaspa"
$ws.Cells.Item(45, 4).Value = 0.04892129
$ws.Cells.Item(45, 5).Value = 1038032171
$ws.Cells.Item(45, 6).Value = 12762572
$ws.Cells.Item(45, 7).Value = -1.44755

$ws.Cells.Item(46, 2).Value = "BSV"
$ws.Cells.Item(46, 3).Value = "Bitcoin SV"
$ws.Cells.Item(46, 4).Value = 51.94
$ws.Cells.Item(46, 5).Value = 1014313569
$ws.Cells.Item(46, 6).Value = 147757498
$ws.Cells.Item(46, 7).Value = -6.68381

$ws.Cells.Item(47, 2).Value = "NEAR"
$ws.Cells.Item(47, 3).Value = "NEAR Protocol"
$ws.Cells.Item(47, 4).Value = 1.03
$ws.Cells.Item(47, 5).Value = 1013833104
$ws.Cells.Item(47, 6).Value = 52755402
$ws.Cells.Item(47, 7).Value = 0.73341

$ws.Cells.Item(48, 2).Value = "RETH"
$ws.Cells.Item(48, 3).Value = "Rocket Pool ETH"
$ws.Cells.Item(48, 4).Value = 1768.61
$ws.Cells.Item(48, 5).Value = 948879933
$ws.Cells.Item(48, 6).Value = 16197010
$ws.Cells.Item(48, 7).Value = 1.3773

$ws.Cells.Item(49, 2).Value = "STX"
$ws.Cells.Item(49, 3).Value = "Stacks"
$ws.Cells.Item(49, 4).Value = 0.620717
$ws.Cells.Item(49, 5).Value = 877865760
$ws.Cells.Item(49, 6).Value = 48132465
$ws.Cells.Item(49, 7).Value = -4.11212

$ws.Cells.Item(50, 4).Value = 5.38
$ws.Cells.Item(50, 5).Value = 775002566
$ws.Cells.Item(50, 6).Value = 7474904
$ws.Cells.Item(50, 7).Value = 0.57728

$ws.Cells.Item(51, 2).Value = "GRT"
$ws.Cells.Item(51, 3).Value = "The Graph"
$ws.Cells.Item(51, 4).Value = 0.08258799999999999
$ws.Cells.Item(51, 5).Value = 764581909
$ws.Cells.Item(51, 6).Value = 32318616
$ws.Cells.Item(51, 7).Value = 1.02293

